# Scheduled runner update: refresh market-price derived columns (H:N)
# across the per-job Sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

function Set-Row($ws, $row, $values) {
    foreach ($col in $values.Keys) {
        $ws.Cells.Item($row, $col).Value = $values[$col]
    }
}

# Column map: H=8 I=9 J=10 K=11 L=12 M=13 N=14

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")

Set-Row $ws 4   @{9=125; 10=15; 11=125; 12=15; 13=-11; 14=-243}
Set-Row $ws 70  @{8=1904.8667; 9=1537.4; 10=2088.6; 11=4612.200000000001; 12=6265.799999999999; 13=-4342.200000000001; 14=-6805.799999999999}
Set-Row $ws 73  @{8=1904.8667; 9=1537.4; 10=2088.6; 11=4612.200000000001; 12=6265.799999999999; 13=-3676.200000000001; 14=-8137.799999999999}
Set-Row $ws 74  @{8=9766.6; 9=8083.1665; 10=10888.889; 11=8083.1665; 12=10888.889; 13=-7147.1665; 14=-12760.889}
Set-Row $ws 77  @{8=9766.6; 9=8083.1665; 10=10888.889; 11=40415.8325; 12=54444.44499999999; 13=-35735.8325; 14=-63804.44499999999}
Set-Row $ws 98  @{8=1534.2858; 9=1462.091; 10=1799; 11=1462.091; 12=1799; 13=35.90900000000011; 14=-4795}
Set-Row $ws 122 @{8=1534.2858; 9=1462.091; 10=1799; 11=4386.272999999999; 12=5397; 13=-1936.272999999999; 14=-10297}
Set-Row $ws 138 @{8=2784.873; 9=1449.3077; 10=3723.3784; 11=4347.9231; 12=11170.1352; 13=792.0769; 14=-21450.1352}
Set-Row $ws 141 @{8=5902.696; 9=5473.125; 10=6884.5713; 11=16419.375; 12=20653.7139; 13=-11239.375; 14=-31013.7139}

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")

Set-Row $ws 32  @{8=20411848; 9=21742826; 11=21742826; 13=-21742539}
Set-Row $ws 132 @{8=1941.5; 9=1941.5; 11=5824.5; 13=-3294.5}

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")

Set-Row $ws 105 @{8=2833.1667; 10=1000; 12=1000; 14=-4494}
Set-Row $ws 132 @{8=3719.1765; 9=3125.1; 10=4567.857; 11=9375.299999999999; 12=13703.571; 13=-6845.299999999999; 14=-18763.571}

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")

Set-Row $ws 8   @{8=1857.25; 9=1857.25; 11=5571.75; 13=-5432.75}

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")

Set-Row $ws 5   @{8=10000; 9=10000; 11=10000; 13=-9888}
Set-Row $ws 15  @{8=0; 10=0; 12=0}
$ws.Cells.Item(15, 14).ClearContents()
Set-Row $ws 81  @{8=0; 10=0; 12=0}
$ws.Cells.Item(81, 14).ClearContents()
Set-Row $ws 84  @{8=0; 10=0; 12=0}
$ws.Cells.Item(84, 14).ClearContents()
Set-Row $ws 102 @{8=3698.4814; 9=2584.5908; 11=2584.5908; 13=-962.5907999999999}
Set-Row $ws 113 @{8=10339.9; 9=4380; 10=16299.8; 11=4380; 12=16299.8; 13=-2210; 14=-20639.8}
Set-Row $ws 132 @{8=6163.2; 9=7231.727; 10=3224.75; 11=21695.181; 12=9674.25; 13=-19165.181; 14=-14734.25}

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")

Set-Row $ws 22  @{8=1833; 9=1250; 10=2124.5; 11=1250; 12=2124.5; 13=-955; 14=-2714.5}
Set-Row $ws 27  @{8=1833; 9=1250; 10=2124.5; 11=1250; 12=2124.5; 13=-1143; 14=-2338.5}
Set-Row $ws 40  @{8=14814.444; 9=34669.332; 11=34669.332; 13=-34533.332}
Set-Row $ws 132 @{8=2220.9768; 9=1871.6061; 10=3373.9; 11=5614.8183; 12=10121.7; 13=-3084.8183; 14=-15181.7}

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")

Set-Row $ws 29  @{8=8003; 9=9504.5; 10=5000; 11=9504.5; 12=5000; 13=-9214.5; 14=-5580}
Set-Row $ws 128 @{8=58635.816; 10=58635.816; 12=58635.816; 14=-68595.81599999999}
Set-Row $ws 132 @{8=1524.6923; 9=1524.6923; 11=4574.0769; 13=-2044.0769}
